# "Assembling And Dissasembling" — revision edits.
#
# Word COM quirk observed in this runtime: any Range.Text mutation (via
# Find/Replace or direct assignment) re-normalises the paragraph by
# merging adjacent runs that share identical formatting. So every textual
# edit is performed FIRST, while text is still naturally coalesced, and
# only afterwards do we carve out the extra run boundaries the target
# markup needs. Run boundaries are produced by dropping a temporary
# bookmark at the split point and immediately deleting it again — Word
# splits the host run to carry the bookmark, and removing the bookmark
# afterwards leaves the (now permanent) run boundary with no formatting
# residue behind.

$d = $word.ActiveDocument

function Find-Range([string]$text) {
    $r = $d.Content
    $null = $r.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $r
}

function Split-At([string]$anchorText, [int]$collapseDirection) {
    # collapseDirection: 0 = wdCollapseEnd, 1 = wdCollapseStart
    $r = Find-Range $anchorText
    $r.Collapse($collapseDirection)
    $bm = "TmpSplit_" + [guid]::NewGuid().ToString("N")
    $d.Bookmarks.Add($bm, $r)
    $d.Bookmarks($bm).Delete()
}

# =====================================================================
# Phase 1 — text content edits
# =====================================================================

# 1. Fix the "Spceifically" typo.
(Find-Range "Spceifically").Text = "Specifically"

# 2. Parenthesise the clarifying clause about disassembly.
(Find-Range "assembly, which turns human-readable instructions into binary, and").Text = `
    "assembly (which turns human-readable instructions into binary) and"

# 3. Add the word "and " before the second half of the asm-embedding sentence.
(Find-Range "to do that we need to add a few extensions").Text = `
    "and to do that we need to add a few extensions"

# 4. Drop the old "_GoBack" bookmark (its paragraph becomes empty).
$d.Bookmarks("_GoBack").Delete()

# =====================================================================
# Phase 2 — run-boundary splits (+ re-placing the "_GoBack" bookmark)
# =====================================================================

# "So, what's ... course!  " | "Specifically" | ", before we write ..."
Split-At "Go off on another tangent of course!  " 0
Split-At "Specifically" 0

# "... an inverse operation to assembly " | "(" | "which turns ... binary" | ")" | " and that is the process ..."
Split-At "an inverse operation to assembly " 0
Split-At "(which turns human-readable instructions into binary" 1
Split-At "(which turns human-readable instructions into binary" 0
Split-At ") and that is the process" 0

# "... our .tim files, " | "and " | [bookmark] | "to do that we need ..."
Split-At "our .tim files, " 0

$r = Find-Range "and to do that we need to add a few extensions"
$r.Collapse(1)          # start of "and to do that ..."
$r.MoveEnd(1, 4)         # extend across exactly "and " (4 characters)
$d.Bookmarks.Add("_GoBack", $r)
